$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-19 16:19:07"
$wsZhCn.Range("E4").Value = "2016-03-19 16:19:07"
$wsZhCn.Range("H3").Value = "2016-03-19 16:19:33"
$wsZhCn.Range("H4").Value = "2016-03-19 16:19:33"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-19 16:19:11"
$wsDeDe.Range("E4").Value = "2016-03-19 16:19:11"
$wsDeDe.Range("H3").Value = "2016-03-19 16:19:38"
$wsDeDe.Range("H4").Value = "2016-03-19 16:19:38"
